$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "git" workflow column in P, mirroring the existing L column
# (git add . / git commit / git push) with the updated steps used on the
# user's home page / dashboard.
$ws.Range("P10").Value = "git status"
$ws.Range("P11").Value = "git add ."
$ws.Range("P12").Value = "git commit -m """""
$ws.Range("P13").Value = "git push origin main"

# Match the selection left behind by the author after making the edit.
$ws.Range("P14").Select()
